$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1. "target" list (column A): insert a new entry "tn.5250" at A26, pushing
#    the existing entries (web, webalert, webcookie, ws, ws.async, xml) down
#    one row (A26:A31 -> A27:A32).
# ---------------------------------------------------------------------------
$targetTail = @("web", "webalert", "webcookie", "ws", "ws.async", "xml")
for ($i = $targetTail.Length - 1; $i -ge 0; $i--) {
    $row = 26 + $i
    $ws.Range("A" + ($row + 1)).Value2 = $targetTail[$i]
}
$ws.Range("A26").Value2 = "tn.5250"

# ---------------------------------------------------------------------------
# 2. "image" list (column K): rename the existing colorbit entry, and insert
#    a new "ocr(image,saveVar)" entry before resize/saveDiff, pushing those
#    two down by one row (K6:K7 -> K7:K8).
# ---------------------------------------------------------------------------
$ws.Range("K8").Value2 = "saveDiff(var,baseline,actual)"
$ws.Range("K7").Value2 = "resize(image,width,height,saveTo)"
$ws.Range("K6").Value2 = "ocr(image,saveVar)"
$ws.Range("K2").Value2 = "colorbit(image,bit,saveTo)"

# ---------------------------------------------------------------------------
# 3. New "tn.5250" list: insert a whole new column before Z (shifting the
#    former Z:AE columns to AA:AF), then populate the new column Z with the
#    command header and its five functions.
# ---------------------------------------------------------------------------
$ws.Columns("Z").Insert()

$ws.Range("Z1").Value2 = "tn.5250"
$ws.Range("Z2").Value2 = "close(profile)"
$ws.Range("Z3").Value2 = "open(profile)"
$ws.Range("Z4").Value2 = "saveText(profile,var)"
$ws.Range("Z5").Value2 = "typeKeys(profile,keystrokes)"
$ws.Range("Z6").Value2 = "updateScreenFields(profile)"

# ---------------------------------------------------------------------------
# 4. Update the workbook-level defined names to reflect the new ranges, and
#    register the new "tn.5250" named range.
# ---------------------------------------------------------------------------
$wb.Names.Item("image").RefersTo = "='#system'!`$K`$2:`$K`$8"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$144"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"
$wb.Names.Add("tn.5250", "='#system'!`$Z`$2:`$Z`$6")
